$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write the new lookup-definition translation rows (531-540).
# Text cells are written in the exact order their underlying shared-string
# values were first introduced in the authored workbook, so that the
# resulting xl/sharedStrings.xml table gets the same new entries appended
# in the same order (indices 1928-1960).
$ws.Range("A531").Value = "TitleSingular"
$ws.Range("F531").Value = "الاسم المفرد"
$ws.Range("F532").Value = "اسم الجمع"
$ws.Range("E531").Value = "Singular Name"
$ws.Range("E532").Value = "Plural Name"
$ws.Range("A532").Value = "PluralName"
$ws.Range("A533").Value = "MainMenuIcon"
$ws.Range("A534").Value = "MainMenuSection"
$ws.Range("A535").Value = "MainMenuSortKey"
$ws.Range("E533").Value = "Icon"
$ws.Range("F535").Value = "الترتيب في القائمة الرئيسية"
$ws.Range("F534").Value = "القسم من القائمة الرئيسية"
$ws.Range("F533").Value = "الأيقونة"
$ws.Range("E534").Value = "Menu Section"
$ws.Range("E535").Value = "Menu Sort Key"
$ws.Range("A536").Value = "UpdateState"
$ws.Range("E536").Value = "Update State"
$ws.Range("F536").Value = "تعديل الحالة"
$ws.Range("G531").Value = "奇异名称"
$ws.Range("G532").Value = "复数名称"
$ws.Range("G533").Value = "图标"
$ws.Range("G534").Value = "菜单部分"
$ws.Range("G535").Value = "菜单排序关键字"
$ws.Range("G536").Value = "更新状态"
$ws.Range("A537").Value = "Definition_State"
$ws.Range("A538").Value = "Definition_State_Draft"
$ws.Range("A539").Value = "Definition_State_Deployed"
$ws.Range("A540").Value = "Definition_State_Archived"
$ws.Range("E539").Value = "Deployed"
$ws.Range("E540").Value = "Archived"
$ws.Range("F539").Value = "مفعل"
$ws.Range("G539").Value = "部署"
$ws.Range("G540").Value = "存档"

# Cells below reuse translation strings that already existed elsewhere in
# the workbook (State / Draft / etc.), so they don't introduce new shared
# strings.
$ws.Range("E537").Value = "State"
$ws.Range("F537").Value = "الحالة"
$ws.Range("G537").Value = "州"
$ws.Range("E538").Value = "Draft"
$ws.Range("F538").Value = "جديد"
$ws.Range("G538").Value = "草案"
$ws.Range("F540").Value = "مؤرشف"

# --- Step 2: Server? / Client? columns are TRUE for every new row.
$ws.Range("B531:C531").Value = $true
$ws.Range("B532:C532").Value = $true
$ws.Range("B533:C533").Value = $true
$ws.Range("B534:C534").Value = $true
$ws.Range("B535:C535").Value = $true
$ws.Range("B536:C536").Value = $true
$ws.Range("B537:C537").Value = $true
$ws.Range("B538:C538").Value = $true
$ws.Range("B539:C539").Value = $true
$ws.Range("B540:C540").Value = $true

# --- Step 3: restore the on-screen selection to roughly where the
# author's cursor ended up after adding the new rows.
$ws.Range("E529").Select()
